$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the header cells belonging to the two ResNet50 sub-tables (fc7 / fc8)
# that are about to be removed, so the merge info doesn't linger after the
# table + data are gone.
$ws.Range("D46:E46").UnMerge()
$ws.Range("G46:H46").UnMerge()
$ws.Range("D47:E47").UnMerge()
$ws.Range("G47:H47").UnMerge()

# Delete the two extra ResNet50 sub-tables (fc7 / fc8); keep "Tabela37" (resnet50)
$ws.ListObjects.Item("Tabela378").Delete()
$ws.ListObjects.Item("Tabela379").Delete()

# Clear the now-orphaned header labels for those two sub-tables (rows 46-47)
$ws.Range("D46:H47").Clear()

# Record the new accuracy reading for the 80% split
$ws.Range("B57").NumberFormat = "@"
$ws.Range("B57").Value = "0.9033"

# Update the visible scroll position / selection to match the author's view
$ws.Range("L54").Select()
